$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents (and formatting) of the "Approved/Rejected" (I) and
# "ReasonToReject" (J) data cells for rows 2 through 30. The column
# headers in row 1 (I1/J1) are left untouched.
$ws.Range("I2:J30").Clear()

# Update the active selection to match the saved state recorded in the
# edited workbook.
$ws.Range("M11").Select()
